$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "Save" header in H1, reusing the same formatting as the other headers (e.g. G1)
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# Fill in the new "Save" column values (H2:H10)
$values = @(1, 0, 0, 1, 0, 1, 1, 1, 0)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $values[$i]
}
